$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# New quest-board rows appended after the existing data (rows 184-196)
# Columns: B = running id (continues the B-column counter as a formula),
#          C = English string, D = Russian translation, E = "converted" string
# ------------------------------------------------------------------

$data = @(
  @("Search requested!",        "Запрошен поиск!",            "Èàðñïšåî ðïéòë!"),
  @("Help! Find Celebi!",       "На помощь! Найдите Селеби!",  "Îà ðïíïþû! Îàêäéóå Òåìåáé!"),
  @("Reward offered!",          "Предлагается награда!",       "Ðñåäìàãàåóòÿ îàãñàäà!"),
  @("Search requested!",        "Запрошен поиск!",             "Èàðñïšåî ðïéòë!"),
  @("Search for Rotom!",        "Найдите Ротома!",             "Îàêäéóå Ñïóïíà!"),
  @("Reward offered!",          "Предлагается награда!",       "Ðñåäìàãàåóòÿ îàãñàäà!"),
  @("Search requested!",        "Запрошен поиск!",             "Èàðñïšåî ðïéòë!"),
  @("Search for Mewtwo!",       "Найдите Мьюту!",              "Îàêäéóå Íûýóô!"),
  @("Reward offered!",          "Предлагается награда!",       "Ðñåäìàãàåóòÿ îàãñàäà!"),
  @("Search requested!",        "Запрошен поиск!",             "Èàðñïšåî ðïéòë!"),
  @("Find Lucario!",            "Найдите Лукарио!",            "Îàêäéóå Ìôëàñéï!"),
  @("Reward offered!",          "Предлагается награда!",       "Ðñåäìàãàåóòÿ îàãñàäà!")
)

$firstRow = 184
$lastDataRow = $firstRow + $data.Count - 1   # 195
$blankRow = $lastDataRow + 1                  # 196

# First row: plain starting number (continues the existing id sequence)
$ws.Range("B$firstRow").Value = 17782

for ($i = 0; $i -lt $data.Count; $i++) {
  $r = $firstRow + $i
  if ($r -gt $firstRow) {
    $prev = $r - 1
    $ws.Range("B$r").Formula = "=B$prev+1"
  }
  $row = $data[$i]
  $ws.Range("C$r").Value = $row[0]
  $ws.Range("D$r").Value = $row[1]
  $ws.Range("E$r").Value = $row[2]
}

# Trailing formatted-but-empty cell, matching the style used in the column
$ws.Range("B2").Copy()
$ws.Range("B$lastDataRow").PasteSpecial(-4122)
for ($r = $firstRow + 1; $r -le $lastDataRow; $r++) {
  $ws.Range("B2").Copy()
  $ws.Range("B$r").PasteSpecial(-4122)
  $prev = $r - 1
  $ws.Range("B$r").Formula = "=B$prev+1"
}
$ws.Range("B2").Copy()
$ws.Range("B$blankRow").PasteSpecial(-4122)
$ws.Range("B$blankRow").ClearContents()

# ------------------------------------------------------------------
# Column E width adjustment
# ------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 30.8333

# ------------------------------------------------------------------
# View state: scroll position & selection, matching where the user
# ended up after adding the rows above
# ------------------------------------------------------------------
$ws.Application.Goto($ws.Range("C196"), $true)
$ws.Range("C196").Select()
$excel.ActiveWindow.ScrollRow = 184
$excel.ActiveWindow.ScrollColumn = 2
